$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Scratch cell used to stage text values so Excel/IronCalc's "smart" numeric
# parsing never kicks in: we force General.NumberFormat="@" (text) on the
# scratch cell, set the literal string there, copy it, then Paste Special
# (values only) into the real target cell. PasteSpecial carries over the
# *text* cell-type as-is (no re-parsing), and the target cell's own style
# (index 0 / no explicit s= attribute) is left untouched. Scratch cell is
# fully Clear()-ed (value + formatting) after every use so it never leaves a
# trace - A1 was blank before the edit and is blank again after.
$scratch = $ws.Range("A1")
function Set-TextValue([string]$cellRef, [string]$text) {
    $scratch.NumberFormat = "@"
    $scratch.Value = $text
    $scratch.Copy()
    $ws.Range($cellRef).PasteSpecial(-4163)
    $scratch.Clear()
}

Set-TextValue "D2" "68.514.44"
Set-TextValue "E2" "  -0.07%  "
Set-TextValue "D3" "2.455.55"
Set-TextValue "E3" "  -0.14%  "
Set-TextValue "E4" "  +0.04%  "
Set-TextValue "D5" "557.20"
Set-TextValue "E5" "  -1.06%  "
Set-TextValue "D6" "160.46"
Set-TextValue "E6" "  -2.13%  "
Set-TextValue "E7" "  +0.08%  "
Set-TextValue "E8" "  +0.02%  "
Set-TextValue "E9" "  -1.49%  "
Set-TextValue "E10" "  +0.58%  "
Set-TextValue "D11" "4.84"
Set-TextValue "E11" "  +0.44%  "
Set-TextValue "E12" "  -3.30%  "
Set-TextValue "D13" "68.433.16"
Set-TextValue "E13" "  -0.03%  "
Set-TextValue "E14" "  -2.65%  "
Set-TextValue "D15" "23.29"
Set-TextValue "E15" "  -1.38%  "
Set-TextValue "D16" "10.58"
Set-TextValue "E16" "  -3.80%  "
Set-TextValue "D17" "332.93"
Set-TextValue "E17" "  -3.04%  "
Set-TextValue "D18" "6.88"
Set-TextValue "E18" "  -3.95%  "
Set-TextValue "D19" "3.76"
Set-TextValue "E19" "  -1.62%  "
Set-TextValue "E20" "  +0.07%  "
Set-TextValue "D21" "1.86"
Set-TextValue "E21" "  -0.98%  "
Set-TextValue "D22" "66.32"
Set-TextValue "E22" "  -2.63%  "
Set-TextValue "D23" "3.61"
Set-TextValue "E23" "  -3.59%  "
Set-TextValue "D24" "8.10"
Set-TextValue "E24" "  -1.71%  "
Set-TextValue "D25" "0.0₃0809"
Set-TextValue "E25" "  -3.81%  "
Set-TextValue "D26" "7.15"
Set-TextValue "E26" "  -2.30%  "
Set-TextValue "E27" "  +0.02%  "
Set-TextValue "D28" "424.96"
Set-TextValue "E28" "  -2.53%  "
Set-TextValue "E29" "  -4.41%  "
Set-TextValue "D30" "1.60"
Set-TextValue "E30" "  -4.67%  "
Set-TextValue "D31" "157.68"
Set-TextValue "E31" "  +0.46%  "
Set-TextValue "E32" "  -0.14%  "
Set-TextValue "E33" "  -0.03%  "
Set-TextValue "E34" "  -1.28%  "
Set-TextValue "D35" "17.69"
Set-TextValue "E35" "  -1.33%  "
Set-TextValue "D36" "0.299"
Set-TextValue "E36" "  -2.87%  "
Set-TextValue "D37" "4.37"
Set-TextValue "E37" "  -2.79%  "
Set-TextValue "E38" "  -5.43%  "
Set-TextValue "E39" "  -3.02%  "
Set-TextValue "E40" "  -2.42%  "
Set-TextValue "D41" "3.32"
Set-TextValue "E41" "  -1.73%  "
Set-TextValue "D42" "128.53"
Set-TextValue "E42" "  -4.65%  "
Set-TextValue "E43" "  -0.25%  "
Set-TextValue "D44" "0.479"
Set-TextValue "E44" "  -1.76%  "
Set-TextValue "E45" "  -0.67%  "
Set-TextValue "E46" "  -0.70%  "
Set-TextValue "E47" "  +0.22%  "
Set-TextValue "D48" "1.37"
Set-TextValue "E48" "  -4.59%  "
Set-TextValue "D49" "4.88"
Set-TextValue "E49" "  -9.55%  "
Set-TextValue "D50" "16.66"
Set-TextValue "E50" "  -5.90%  "
Set-TextValue "D51" "0.0₆0204"
Set-TextValue "E51" "  -1.16%  "

$wb.Application.CutCopyMode = $false

